$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark the newly-completed "x" cells in column D ---
$ws.Range("D45").Value = "x"
$ws.Range("D46").Value = "x"
$ws.Range("D76").Value = "x"
$ws.Range("D77").Value = "x"
$ws.Range("D78").Value = "x"
$ws.Range("D79").Value = "x"
$ws.Range("D80").Value = "x"
$ws.Range("D81").Value = "x"

# --- Re-apply Print Titles (row 1) on this sheet; LibreOffice/Excel keep
#     appending a fresh _xlnm.Print_Titles_N defined name each time the
#     Page Setup "Rows to repeat at top" is (re)confirmed, so mint the next
#     one in the sequence, scoped to this worksheet. ---
$ws.Names.Add("_xlnm.Print_Titles_0_0_0_0_0", "='Casos de Uso'!`$1:`$1")

# --- Scroll/selection state: move the viewport back to the top and put
#     the active cell on D49 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("D49").Select() | Out-Null
